$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("14369201", "2025-08-15", "Laurent Lokoli", "Dan Added", "Gana Dan Added", 2.1),
    @("14360037", "2025-08-15", "Nerman Fatic", "Zdenek Kolar", "Gana Zdenek Kolar", 3.25),
    @("14370728", "2025-08-14", "James McCabe", "Lukas Klein", "Gana James McCabe", 2.75)
)

$startRow = 201
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # event_id (A) must stay text (e.g. "14369201"), not be coerced to a number
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 1).Style = "Normal"

    # fecha (B) must stay text (e.g. "2025-08-15"), not be coerced to a date
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]

    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = ""
}
